# Add two new "tutorial" (practice) blocks of questions to the question bank.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A = Paragraph_id, B = Question_id, C = Question, D = Answers, E = Correct_answer

# --- Block 1: "practice1" (rows 14-16) ---
$ws.Range("A14:A16").Value = "practice1"

# --- Block 2: "practice2" (rows 17-19) ---
$ws.Range("A17:A19").Value = "practice2"

$ws.Range("C14").Value = "1.Where was the city of Oceana located?"
$ws.Range("C15").Value = "2. What did Jenna and Max find inside the ancient shipwreck?"
$ws.Range("D15").Value = "a) A treasure chest of gold`nb) Old maps and strange artifacts`nc) A secret message`nd) A new energy source"
$ws.Range("C16").Value = "3. How did the discovery of the shipwrecks benefit Oceana?"
$ws.Range("D14").Value = "a) In the Arctic Ocean`nb) On the moon`nc) In the Pacific Ocean`nd) On a mountaintop"
$ws.Range("D16").Value = "a) It made Jenna famous`nb) It provided valuable resources`nc) It led to the construction of another city`nd) It stopped the sea levels from rising"

$ws.Range("C17").Value = "1.What was Verdantia known for?"
$ws.Range("D17").Value = "a) Its giant trees and fertile soil`nb) Its vast deserts`nc) Its advanced technology`nd) Its floating cities"
$ws.Range("C18").Value = "2. What did Liam discover in the forest?"
$ws.Range("D18").Value = "a) A hidden treasure`nb) A new species of animal`nc) A sick tree showing signs of decay`nd) A secret cave"
$ws.Range("C19").Value = "3. How did the Galactic Gardeners save the trees from the fungus?"
$ws.Range("D19").Value = "a) By cutting down the infected trees`nb) By relocating the trees`nc) By developing a special serum`nd) By using robots to clean the trees"

# Fill in the Question_id (B) and Correct_answer (E) columns, reusing the existing
# shared-string values (Q1/Q2/Q3 and a/b/c/d) already present in the workbook.
$ws.Range("B14").Value = "Q1"
$ws.Range("E14").Value = "c"
$ws.Range("B15").Value = "Q2"
$ws.Range("E15").Value = "b"
$ws.Range("B16").Value = "Q3"
$ws.Range("E16").Value = "b"
$ws.Range("B17").Value = "Q1"
$ws.Range("E17").Value = "a"
$ws.Range("B18").Value = "Q2"
$ws.Range("E18").Value = "c"
$ws.Range("B19").Value = "Q3"
$ws.Range("E19").Value = "c"

# Wrap the answer-choices column and size the new rows the same as the existing
# question rows (4 wrapped answer lines at the default 14.4pt row height).
$ws.Range("D14:D19").WrapText = $true
for ($n = 14; $n -le 19; $n++) {
    $ws.Rows.Item($n).RowHeight = 57.6
}

# Update the view: scroll so row 15 is the top visible row, and select G17 (matches author's last-saved view state)
$ws.Range("G17").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 1
